# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strikeout) values replacing the previous Strike# counts, row by row (row 2 .. row 28)
$kValues = @{
    2  = 1
    3  = 6
    4  = 4
    5  = 4
    6  = 3
    7  = 7
    8  = 6
    9  = 6
    10 = 2
    11 = 7
    12 = 9
    13 = 3
    14 = 1
    15 = 3
    16 = 4
    17 = 6
    18 = 4
    19 = 2
    20 = 3
    21 = 4
    22 = 1
    23 = 6
    24 = 5
    25 = 3
    26 = 2
    27 = 4
    28 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
